# Refs #5722. Refs #5723.
# Column re-alignment corrections: insert a new "Hire Date" column before
# the existing "Benefit Begin Date" column (column H) on the active sheet,
# shifting all subsequent columns one position to the right, and leave the
# selection on the newly inserted column's data cell (H2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at H; everything from H onward (including the
# old H column) shifts right by one.
$ws.Columns("H").Insert() | Out-Null

# Label the new column's header.
$ws.Range("H1").Value = "Hire Date"

# Leave the active selection on the new column in the data row, matching
# the saved view state.
$ws.Range("H2").Select() | Out-Null
